$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New submission link shared by the two students whose TRAB 6 grade / GitHub
# link are being updated in this pass of grading.
$newLink = "https://github.com/IgorSantos26/Estoque-sobrevivencia"

# Row 6 - IGOR OLIVEIRA SANTOS: TRAB 6 grade corrected from 0 to 3, and the
# App Cloud / GitHub submission is now recognised (mark the row as graded
# with the light "reviewed" shading, same as other graded rows).
$ws.Range("G6").Value = 3
$ws.Range("I6").Value = $newLink
$ws.Range("I6").Font.Name = "Times New Roman"
$ws.Range("I6").Font.Size = 11
$ws.Range("A6").Interior.ThemeColor = 2

# Row 13 - PAULA LETICIA QUEIROZ DOS ANJOS: same TRAB 6 correction and link.
$ws.Range("G13").Value = 3
$ws.Range("I13").Value = $newLink
$ws.Range("A13").Interior.ThemeColor = 2

$wb.Application.Calculate()

# Restore the view to the top of the sheet and leave the selection on A11.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A11").Select()
